$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the leading "Unidade acadêmica / curso / habilitação" / index column
# (column A). Everything to its right (Raca, Feminino, % Feminino, ...)
# shifts one column to the left.
$ws.Columns("A").Delete()
